$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.2176170349121094
$ws.Range("E2").Value = 74.54980530067041
$ws.Range("F2").Value = 0.003124110796206647
$ws.Range("G2").Value = 0.002501771280629178
$ws.Range("H2").Value = 0.002371583497738418
$ws.Range("I2").Value = 0.002190485478716327
$ws.Range("J2").Value = 0.002190485478716327
$ws.Range("K2").Value = 0.002082988396314586
$ws.Range("L2").Value = 0.002065289485462057
$ws.Range("M2").Value = 0.002065289485462057
$ws.Range("N2").Value = 0.002065289485462057
$ws.Range("O2").Value = 0.001943652715824293
$ws.Range("P2").Value = 0.001943652715824293
$ws.Range("Q2").Value = 0.001885223026931721
$ws.Range("R2").Value = 0.001773839022968643
$ws.Range("S2").Value = 0.001637779492893455
$ws.Range("T2").Value = 0.001615034114461364
$ws.Range("U2").Value = 0.001559535802326531
$ws.Range("V2").Value = 0.001530722176329056
$ws.Range("W2").Value = 0.001489872629669288
$ws.Range("X2").Value = 0.001476237711573956
$ws.Range("Y2").Value = 0.001453212578960437

# Row 3
$ws.Range("C3").Value = 0.2034707069396973
$ws.Range("E3").Value = 77.47709094679522
$ws.Range("F3").Value = 0.003119161324487111
$ws.Range("G3").Value = 0.002614634180354116
$ws.Range("H3").Value = 0.002207262044267145
$ws.Range("I3").Value = 0.00215950341019966
$ws.Range("J3").Value = 0.002123073436809416
$ws.Range("K3").Value = 0.001812903825000752
$ws.Range("L3").Value = 0.001812903825000752
$ws.Range("M3").Value = 0.001812903825000752
$ws.Range("N3").Value = 0.001812903825000752
$ws.Range("O3").Value = 0.001781644122554362
$ws.Range("P3").Value = 0.001759730198981889
$ws.Range("Q3").Value = 0.001682123083133228
$ws.Range("R3").Value = 0.001665349034886455
$ws.Range("S3").Value = 0.001630174913366858
$ws.Range("T3").Value = 0.001559890766084226
$ws.Range("U3").Value = 0.001559890766084226
$ws.Range("V3").Value = 0.001559890766084226
$ws.Range("W3").Value = 0.001540701445962988
$ws.Range("X3").Value = 0.001524158263710297
$ws.Range("Y3").Value = 0.001510274677325443

# Row 4
$ws.Range("C4").Value = 0.2069470882415771
$ws.Range("E4").Value = 77.94497918227353
$ws.Range("F4").Value = 0.003154521551819769
$ws.Range("G4").Value = 0.002546079571777659
$ws.Range("H4").Value = 0.002515459946799986
$ws.Range("I4").Value = 0.002277736284254095
$ws.Range("J4").Value = 0.002111207500557925
$ws.Range("K4").Value = 0.002111207500557925
$ws.Range("L4").Value = 0.002031484973807214
$ws.Range("M4").Value = 0.001923826695350917
$ws.Range("N4").Value = 0.001905164871156325
$ws.Range("O4").Value = 0.001880209174435409
$ws.Range("P4").Value = 0.001753232597861922
$ws.Range("Q4").Value = 0.001753232597861922
$ws.Range("R4").Value = 0.001658618783071838
$ws.Range("S4").Value = 0.001624584513387974
$ws.Range("T4").Value = 0.001610028858849475
$ws.Range("U4").Value = 0.001610028858849475
$ws.Range("V4").Value = 0.001588022233833175
$ws.Range("W4").Value = 0.001530149233407194
$ws.Range("X4").Value = 0.001530149233407194
$ws.Range("Y4").Value = 0.001519395305697339

# Row 5
$ws.Range("C5").Value = 0.2263250350952148
$ws.Range("E5").Value = 75.54254521028088
$ws.Range("F5").Value = 0.003055602095804392
$ws.Range("G5").Value = 0.002548989490530578
$ws.Range("H5").Value = 0.002429329037792827
$ws.Range("I5").Value = 0.002134699443208073
$ws.Range("J5").Value = 0.001998537836782658
$ws.Range("K5").Value = 0.001998537836782658
$ws.Range("L5").Value = 0.001998537836782658
$ws.Range("M5").Value = 0.001693041315815505
$ws.Range("N5").Value = 0.001693041315815505
$ws.Range("O5").Value = 0.001693041315815505
$ws.Range("P5").Value = 0.001693041315815505
$ws.Range("Q5").Value = 0.001670435300451888
$ws.Range("R5").Value = 0.001617809075837426
$ws.Range("S5").Value = 0.001617809075837426
$ws.Range("T5").Value = 0.001590715522267599
$ws.Range("U5").Value = 0.001556132530445237
$ws.Range("V5").Value = 0.001518729963093305
$ws.Range("W5").Value = 0.001517020222946816
$ws.Range("X5").Value = 0.00148429059636328
$ws.Range("Y5").Value = 0.001472564234118535

# Row 6
$ws.Range("C6").Value = 0.1933751106262207
$ws.Range("E6").Value = 71.59786819917099
$ws.Range("F6").Value = 0.003065921266638414
$ws.Range("G6").Value = 0.002681937696558773
$ws.Range("H6").Value = 0.002418377080786641
$ws.Range("I6").Value = 0.002299187119463108
$ws.Range("J6").Value = 0.002094839229990136
$ws.Range("K6").Value = 0.002094839229990136
$ws.Range("L6").Value = 0.001947447562156314
$ws.Range("M6").Value = 0.001947447562156314
$ws.Range("N6").Value = 0.001772963399753687
$ws.Range("O6").Value = 0.001634807895201869
$ws.Range("P6").Value = 0.001634807895201869
$ws.Range("Q6").Value = 0.001634807895201869
$ws.Range("R6").Value = 0.00149951857504386
$ws.Range("S6").Value = 0.00149951857504386
$ws.Range("T6").Value = 0.00149951857504386
$ws.Range("U6").Value = 0.00149951857504386
$ws.Range("V6").Value = 0.001487233277658041
$ws.Range("W6").Value = 0.001466875879081268
$ws.Range("X6").Value = 0.00141574826245791
$ws.Range("Y6").Value = 0.001395669945402943

# Row 7
$ws.Range("C7").Value = 0.2256424427032471
$ws.Range("E7").Value = 80.14539883980251
$ws.Range("F7").Value = 0.00299220981746574
$ws.Range("G7").Value = 0.002656296134397694
$ws.Range("H7").Value = 0.002497502014548356
$ws.Range("I7").Value = 0.002326162035366549
$ws.Range("J7").Value = 0.002326162035366549
$ws.Range("K7").Value = 0.002118839845021518
$ws.Range("L7").Value = 0.002078752228266076
$ws.Range("M7").Value = 0.00205060481672188
$ws.Range("N7").Value = 0.001950412159824454
$ws.Range("O7").Value = 0.001856462944413412
$ws.Range("P7").Value = 0.001765420526254872
$ws.Range("Q7").Value = 0.001765420526254872
$ws.Range("R7").Value = 0.001680878513559957
$ws.Range("S7").Value = 0.001680878513559957
$ws.Range("T7").Value = 0.001629120250346607
$ws.Range("U7").Value = 0.001609288559614116
$ws.Range("V7").Value = 0.001603098100308502
$ws.Range("W7").Value = 0.001585795928309295
$ws.Range("X7").Value = 0.001562288476409405
$ws.Range("Y7").Value = 0.001562288476409405

# Row 8
$ws.Range("C8").Value = 0.1923003196716309
$ws.Range("E8").Value = 78.53660796625263
$ws.Range("F8").Value = 0.003009124105867431
$ws.Range("G8").Value = 0.002388540814134716
$ws.Range("H8").Value = 0.002388540814134716
$ws.Range("I8").Value = 0.002173010236696299
$ws.Range("J8").Value = 0.002137537169354132
$ws.Range("K8").Value = 0.002056321753901316
$ws.Range("L8").Value = 0.001894050972198484
$ws.Range("M8").Value = 0.001894050972198484
$ws.Range("N8").Value = 0.001871458052208123
$ws.Range("O8").Value = 0.001825144198802825
$ws.Range("P8").Value = 0.001818974870371537
$ws.Range("Q8").Value = 0.00178829615041556
$ws.Range("R8").Value = 0.00178829615041556
$ws.Range("S8").Value = 0.001761116279348639
$ws.Range("T8").Value = 0.001673573969191244
$ws.Range("U8").Value = 0.001545019654426865
$ws.Range("V8").Value = 0.001545019654426865
$ws.Range("W8").Value = 0.001545019654426865
$ws.Range("X8").Value = 0.001544992752320258
$ws.Range("Y8").Value = 0.00153092803053124

# Row 9
$ws.Range("C9").Value = 0.2138075828552246
$ws.Range("E9").Value = 76.60816520298613
$ws.Range("F9").Value = 0.003154521551819769
$ws.Range("G9").Value = 0.002448038540366686
$ws.Range("H9").Value = 0.00210534407484548
$ws.Range("I9").Value = 0.00210534407484548
$ws.Range("J9").Value = 0.001993402778122993
$ws.Range("K9").Value = 0.001993402778122993
$ws.Range("L9").Value = 0.001908239188164655
$ws.Range("M9").Value = 0.001908239188164655
$ws.Range("N9").Value = 0.001892791752908952
$ws.Range("O9").Value = 0.001812644877832302
$ws.Range("P9").Value = 0.001774093044534052
$ws.Range("Q9").Value = 0.001758486381944426
$ws.Range("R9").Value = 0.001628709542550793
$ws.Range("S9").Value = 0.001628709542550793
$ws.Range("T9").Value = 0.001605416923664617
$ws.Range("U9").Value = 0.001572818368508892
$ws.Range("V9").Value = 0.001553027808574173
$ws.Range("W9").Value = 0.001536364285581873
$ws.Range("X9").Value = 0.001519335049124492
$ws.Range("Y9").Value = 0.001493336553664447

# Row 10
$ws.Range("C10").Value = 0.2248187065124512
$ws.Range("E10").Value = 76.05203231859741
$ws.Range("F10").Value = 0.003154521551819769
$ws.Range("G10").Value = 0.002570832371685492
$ws.Range("H10").Value = 0.002252453466438039
$ws.Range("I10").Value = 0.002120590710439064
$ws.Range("J10").Value = 0.001947824114868303
$ws.Range("K10").Value = 0.001947824114868303
$ws.Range("L10").Value = 0.001921417796391459
$ws.Range("M10").Value = 0.001780988954615381
$ws.Range("N10").Value = 0.001780988954615381
$ws.Range("O10").Value = 0.001756091207802523
$ws.Range("P10").Value = 0.001703139096069023
$ws.Range("Q10").Value = 0.001699862849928272
$ws.Range("R10").Value = 0.001624770185480316
$ws.Range("S10").Value = 0.001624770185480316
$ws.Range("T10").Value = 0.001612089265580467
$ws.Range("U10").Value = 0.001612089265580467
$ws.Range("V10").Value = 0.001559986288057526
$ws.Range("W10").Value = 0.001537243009684549
$ws.Range("X10").Value = 0.001499178989643662
$ws.Range("Y10").Value = 0.001482495756697805

# Row 11
$ws.Range("C11").Value = 0.2432901859283447
$ws.Range("E11").Value = 74.22437339502721
$ws.Range("F11").Value = 0.003154521551819769
$ws.Range("G11").Value = 0.002623308732567855
$ws.Range("H11").Value = 0.002307739549751576
$ws.Range("I11").Value = 0.002186604467092835
$ws.Range("J11").Value = 0.002090644353823639
$ws.Range("K11").Value = 0.002090644353823639
$ws.Range("L11").Value = 0.001974019526614637
$ws.Range("M11").Value = 0.001913523822662519
$ws.Range("N11").Value = 0.001906488191992335
$ws.Range("O11").Value = 0.001833016994224541
$ws.Range("P11").Value = 0.001726218639272822
$ws.Range("Q11").Value = 0.001686316111514691
$ws.Range("R11").Value = 0.001630476189793359
$ws.Range("S11").Value = 0.001593644828524943
$ws.Range("T11").Value = 0.001539268788906555
$ws.Range("U11").Value = 0.001539268788906555
$ws.Range("V11").Value = 0.001505352392207171
$ws.Range("W11").Value = 0.001491374766002608
$ws.Range("X11").Value = 0.001452453856487749
$ws.Range("Y11").Value = 0.001446868877096047

